# ============================================================================
# Email-FAQ-en.xlsx KB update
#  - B4 ("file could not be submitted" answer): bold the question, keep the
#    rest (including the existing hyperlink run) as a normal run.
#  - B6 ("view or edit a report" answer): bold only the question (own run),
#    restructure the remaining runs (answer text / existing email link /
#    closing remark) to match the new paragraph breaks.
#  - Column A narrower + wrap + top-align; column B header wraps too.
#  - Row heights grow to fit the re-wrapped text.
#  - Selection moves to A15.
# ============================================================================

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# B4: "Why do I receive an error message "The file could not be
#      submitted because errors were found" ..."
# ---------------------------------------------------------------------
$b4Text = "Why do I receive an error message ""The file could not be submitted because errors were found"" when attempting to submit my monthly report in CTLS? `nPlease refer to the file specifications in the monthly reporting guide as well as the inventory reporting tool. The purpose of this tool is to provide further guidance and assistance in preparing monthly report submission into the Cannabis Tracking and Licensing System (CTLS). This tool does not replace or substitute the .csv template above but rather can be used to populate the .csv file. The reporting guide and the inventory tool are available in the link below: https://www.canada.ca/en/health-canada/services/drugs-medication/cannabis/tracking-system.html"
$ws.Range("B4").Value = $b4Text

$b4Question = $ws.Range("B4").Characters(1, 146)
$b4Question.Font.Bold = $true

$b4Body = $ws.Range("B4").Characters(147, 471)
$b4Body.Font.Bold = $false

$b4Link = $ws.Range("B4").Characters(618, 94)
$b4Link.Font.Bold = $false
$b4Link.Font.Color = 255 * 65536

# ---------------------------------------------------------------------
# B6: "How do I view or edit a report once it has been submitted in
#      CTLS?" ...
# ---------------------------------------------------------------------
$b6Text = "How do I view or edit a report once it has been submitted in CTLS?`nOnce a report has been successfully submitted, it cannot be re-opened for further editing. If a correction to a previously submitted report is required, send a request to hc.ctls-bi-sscdl-ie.sc@canada.ca.`nWe would like to remind you that it is your responsibility to ensure that you retain a copy of your monthly submissions. "
$ws.Range("B6").Value = $b6Text

$b6Question = $ws.Range("B6").Characters(1, 67)
$b6Question.Font.Bold = $true

$b6Answer = $ws.Range("B6").Characters(68, 171)
$b6Answer.Font.Bold = $false

$b6Email = $ws.Range("B6").Characters(239, 32)
$b6Email.Font.Bold = $false
$b6Email.Font.Color = 255 * 65536

$b6Closing = $ws.Range("B6").Characters(271, 123)
$b6Closing.Font.Bold = $false

# ---------------------------------------------------------------------
# Column layout: narrower column A, word-wrap + top vertical align;
# column B header now wraps too.
# ---------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 55.25
$ws.Columns.Item(1).WrapText = $true
$ws.Columns.Item(1).VerticalAlignment = -4160

$ws.Range("B1").WrapText = $true

# ---------------------------------------------------------------------
# Row heights grow now that column A re-wraps with the new width.
# ---------------------------------------------------------------------
$ws.Rows.Item(2).RowHeight = 46.5
$ws.Rows.Item(3).RowHeight = 80.25
$ws.Rows.Item(4).RowHeight = 91.5
$ws.Rows.Item(5).RowHeight = 57.75
$ws.Rows.Item(6).RowHeight = 57.75
$ws.Rows.Item(7).RowHeight = 102.75
$ws.Rows.Item(8).RowHeight = 57.75

# ---------------------------------------------------------------------
# Move the selection, matching the refreshed view.
# ---------------------------------------------------------------------
$ws.Range("A15").Select()
